# Rename the shared "MetaIterable / MetaMap" string to "MetaIterable / MetaMappable"
# wherever it appears on the "meta types" worksheet (std::map / std::multimap /
# std::unordered_map / std::unordered_multimap rows, column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta types")

$oldText = "MetaIterable`nMetaMap"
$newText = "MetaIterable`nMetaMappable"

$targetCells = @("D45", "D46", "D49", "D50")
foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    if ($cell.Value() -eq $oldText) {
        $cell.Value = $newText
    }
}

$ws.Range("A2").Select()
